{"js": "// Update the worksheet date heading and the 25 three-digit-by-one-digit\n// multiplication facts in the table to the new day's values.\n//\n// Every \"old\" string below occurs exactly once in the document, so for each\n// pair we search the body for that exact text and replace the (single) hit\n// with the new text. This is insensitive to the physical location of each\n// run (title paragraph vs. table cell) and to the order pairs are listed in.\nconst replacements = [\n  [\"2024-05-04 Saturday\", \"2024-05-05 Sunday\"],\n  [\"811\u00d78=6488\", \"446\u00d73=1338\"],\n  [\"483\u00d79=4347\", \"875\u00d76=5250\"],\n  [\"383\u00d79=3447\", \"154\u00d76=924\"],\n  [\"895\u00d79=8055\", \"237\u00d78=1896\"],\n  [\"146\u00d76=876\", \"985\u00d72=1970\"],\n  [\"251\u00d78=2008\", \"524\u00d78=4192\"],\n  [\"827\u00d74=3308\", \"553\u00d75=2765\"],\n  [\"923\u00d75=4615\", \"237\u00d79=2133\"],\n  [\"498\u00d75=2490\", \"908\u00d76=5448\"],\n  [\"914\u00d78=7312\", \"868\u00d74=3472\"],\n  [\"823\u00d72=1646\", \"761\u00d72=1522\"],\n  [\"873\u00d75=4365\", \"354\u00d76=2124\"],\n  [\"466\u00d78=3728\", \"178\u00d72=356\"],\n  [\"383\u00d77=2681\", \"867\u00d74=3468\"],\n  [\"841\u00d78=6728\", \"665\u00d77=4655\"],\n  [\"514\u00d74=2056\", \"561\u00d77=3927\"],\n  [\"628\u00d76=3768\", \"372\u00d78=2976\"],\n  [\"816\u00d72=1632\", \"825\u00d73=2475\"],\n  [\"629\u00d74=2516\", \"575\u00d75=2875\"],\n  [\"499\u00d79=4491\", \"450\u00d77=3150\"],\n  [\"221\u00d79=1989\", \"424\u00d73=1272\"],\n  [\"249\u00d73=747\", \"783\u00d79=7047\"],\n  [\"983\u00d77=6881\", \"157\u00d79=1413\"],\n  [\"686\u00d75=3430\", \"945\u00d76=5670\"],\n  [\"541\u00d73=1623\", \"373\u00d75=1865\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found, cannot apply edit: \" + oldText);\n  }\n\n  // Replace every hit (normally exactly one) so the script is still correct\n  // if a value were ever duplicated elsewhere in the document.\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date heading and the 25 multiplication facts in the table\n# with the values from the updated worksheet. Every source string in this\n# document is unique, so Find/Replace (wdReplaceAll, matched on exactly one\n# hit each) is unambiguous and order independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-05-04 Saturday', '2024-05-05 Sunday'),\n    @('811\u00d78=6488', '446\u00d73=1338'),\n    @('483\u00d79=4347', '875\u00d76=5250'),\n    @('383\u00d79=3447', '154\u00d76=924'),\n    @('895\u00d79=8055', '237\u00d78=1896'),\n    @('146\u00d76=876', '985\u00d72=1970'),\n    @('251\u00d78=2008', '524\u00d78=4192'),\n    @('827\u00d74=3308', '553\u00d75=2765'),\n    @('923\u00d75=4615', '237\u00d79=2133'),\n    @('498\u00d75=2490', '908\u00d76=5448'),\n    @('914\u00d78=7312', '868\u00d74=3472'),\n    @('823\u00d72=1646', '761\u00d72=1522'),\n    @('873\u00d75=4365', '354\u00d76=2124'),\n    @('466\u00d78=3728', '178\u00d72=356'),\n    @('383\u00d77=2681', '867\u00d74=3468'),\n    @('841\u00d78=6728', '665\u00d77=4655'),\n    @('514\u00d74=2056', '561\u00d77=3927'),\n    @('628\u00d76=3768', '372\u00d78=2976'),\n    @('816\u00d72=1632', '825\u00d73=2475'),\n    @('629\u00d74=2516', '575\u00d75=2875'),\n    @('499\u00d79=4491', '450\u00d77=3150'),\n    @('221\u00d79=1989', '424\u00d73=1272'),\n    @('249\u00d73=747', '783\u00d79=7047'),\n    @('983\u00d77=6881', '157\u00d79=1413'),\n    @('686\u00d75=3430', '945\u00d76=5670'),\n    @('541\u00d73=1623', '373\u00d75=1865')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
